$d = $word.ActiveDocument

# Map of old text -> new text (all unique, no cascading overlaps)
$replacements = @(
    @("2024-11-29 Friday", "2024-11-30 Saturday"),
    @("60×21=1260", "55×26=1430"),
    @("47×35=1645", "56×77=4312"),
    @("91×57=5187", "38×98=3724"),
    @("25×46=1150", "29×26=754"),
    @("37×73=2701", "39×56=2184"),
    @("12×22=264", "29×17=493"),
    @("97×92=8924", "19×42=798"),
    @("41×52=2132", "15×50=750"),
    @("85×75=6375", "57×90=5130"),
    @("85×84=7140", "68×66=4488"),
    @("81×22=1782", "58×75=4350"),
    @("20×70=1400", "11×77=847"),
    @("29×55=1595", "84×99=8316"),
    @("21×43=903", "75×53=3975"),
    @("35×43=1505", "79×13=1027"),
    @("87×35=3045", "88×89=7832"),
    @("39×16=624", "65×98=6370"),
    @("95×80=7600", "94×15=1410"),
    @("80×64=5120", "38×81=3078"),
    @("70×61=4270", "39×87=3393"),
    @("43×14=602", "34×41=1394"),
    @("52×54=2808", "80×97=7760"),
    @("21×84=1764", "86×20=1720"),
    @("38×61=2318", "18×33=594"),
    @("40×47=1880", "18×84=1512")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
